# DebitCreditNoteForReinsurance - add Reinsured Name / Reinsurer Address columns
# to the CreditNoteDetails sheet, and clear out the now-unused second Debit Note
# row on Sheet1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)            # "Sheet1"            (Debit Note)
$ws2 = $wb.Worksheets.Item(2)            # "CreditNoteDetails" (Credit Note)

# ---------------------------------------------------------------------------
# Sheet1: the second debit-note row (row 3, RI/25-26/GIFT/D20) is removed from
# use - clear its contents but keep the existing cell styling/borders intact.
# ---------------------------------------------------------------------------
$ws1.Range("A3:V3").ClearContents()

# ---------------------------------------------------------------------------
# CreditNoteDetails: insert two new columns - "Reinsured Name" (C) and
# "Reinsurer Address" (E) - between the existing columns.
# ---------------------------------------------------------------------------
$ws2.Columns("C").Insert()
$ws2.Columns("E").Insert()

# Header row
$ws2.Range("C1").Value = "Reinsured Name"
$ws2.Range("E1").Value = "Reinsurer Address"

# Data rows (2-4): every credit note line belongs to the same reinsured
# company and the same reinsurer address.
$reinsuredName    = "Solarelle Insurance Pvt. Ltd"
$reinsurerAddress = "Green City Office Park, Danny Pule Road Lusaka ZM, 10101, Zambia"

$ws2.Range("C2").Value = $reinsuredName
$ws2.Range("C3").Value = $reinsuredName
$ws2.Range("C4").Value = $reinsuredName

$ws2.Range("E2").Value = $reinsurerAddress
$ws2.Range("E3").Value = $reinsurerAddress
$ws2.Range("E4").Value = $reinsurerAddress

# The address column wraps, so the sample rows grow to a two-line height.
$ws2.Range("E2:E4").WrapText = $true
$ws2.Range("C2:C4").Style = "Normal"

$ws2.Rows("2").RowHeight = 28.8
$ws2.Rows("3").RowHeight = 28.8
$ws2.Rows("4").RowHeight = 28.8

# These three figures are no longer carried on this sheet (now blank).
$ws2.Range("I2:I4").ClearContents()
$ws2.Range("K2:K4").ClearContents()
$ws2.Range("N2:N4").ClearContents()

# ---------------------------------------------------------------------------
# Selections / view state, matching where the author left the cursor.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("J2:O4").Select()

$ws1.Activate()
$ws1.Range("L2:V2").Select()
